$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '54.299.55'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +4.59%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.172.38'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.01%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '400.03'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.69%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '109.10'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.66%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.619'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +4.43%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '38.89'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.34%  '
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0881'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.72%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.673.53'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.23%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '19.15'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.81%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.03'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.82%  '
$ws.Range('E16').Value = '  +7.86%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.184.78'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.51'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.13%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '54.319.76'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.33%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.30'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.58%  '
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0989'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.75%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '71.57'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '272.84'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.26'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.30%  '
$ws.Range('E26').Value = '  -2.78%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '27.67'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.169'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.75%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  +2.97%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '11.04'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.24%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0496'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +9.64%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '36.73'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.19%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '50.56'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('E37').Value = '  +6.86%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.83'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +8.73%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.10'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +9.51%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.291'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.65%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.92'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.30%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '17.28'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '130.33'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.04%  '
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.50'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.087.17'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0339'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +5.34%  '
$ws.Range('E51').Value = '  +8.26%  '
